# Generate Report for Handback
#
# The handback transform for file 9fdd8f93-ade8-432b-952b-8081534a9a67
# (row 3 on every sheet) failed, so the status columns flip from
# "Ready for handoff" to "Handback transform failed" and a new
# "Error Detail" (column L) entry is recorded for both locales
# explaining the mismatched handback file name.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 = 9fdd8f93-ade8-432b-952b-8081534a9a67.md
# Column B = zh-cn status, Column C = de-de status
$wsOverview.Cells.Item(3, 2).Value = $newStatus
$wsOverview.Cells.Item(3, 3).Value = $newStatus

# zh-cn sheet: row 3 = 9fdd8f93-ade8-432b-952b-8081534a9a67.md
# Column C = Status, Column L = Error Detail
$wsZhCn.Cells.Item(3, 3).Value = $newStatus
$wsZhCn.Cells.Item(3, 12).Value = "Handback file name: 5d1zbzks.fro is different with handoff file name: 9fdd8f93-ade8-432b-952b-8081534a9a67.f9b6c7b44d3d9179c63b947c71628f1f9ff04949.zh-cn."

# de-de sheet: row 3 = 9fdd8f93-ade8-432b-952b-8081534a9a67.md
# Column C = Status, Column L = Error Detail
$wsDeDe.Cells.Item(3, 3).Value = $newStatus
$wsDeDe.Cells.Item(3, 12).Value = "Handback file name: 5d1zbzks.fro is different with handoff file name: 9fdd8f93-ade8-432b-952b-8081534a9a67.f9b6c7b44d3d9179c63b947c71628f1f9ff04949.de-de."
